# QuizSetupSpreadsheet.xlsx edit
# "Finished Quiz  includes old quiz junk"
#
# Replaces the old placeholder picture paths / example rows with the
# real finished quiz content on Sheet1 ("Table1"):
#   - Row 3 & 4: point the picture path to the new local repo path and
#     re-order which multiple-choice answer text (A/B/C/D) goes with
#     which question.
#   - Row 5: becomes a Text-type question (was SelectOne), new picture path.
#   - Row 6 & 7: fill in two more quiz rows that used to be blank templates
#     (Text & Image / SelectMulti, and Text / True-False questions).
#   - Row heights for rows 3-5 grow from 45 to 75 to fit the longer text.
#   - Selection moves to V5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3 ---------------------------------------------------------------
$ws.Range("C3").Value = "C:\Users\husmith\Documents\GitHub\Action Pack\Quiz\Question1.png"
$ws.Range("K3").Value = "A"
$ws.Range("M3").Value = "B"
$ws.Range("O3").Value = "C"
$ws.Range("Q3").Value = "D"

# --- Row 4 ---------------------------------------------------------------
$ws.Range("C4").Value = "C:\Users\husmith\Documents\GitHub\Action Pack\Quiz\Question2.png"
$ws.Range("K4").Value = "A"
$ws.Range("M4").Value = "B"
$ws.Range("O4").Value = "C"
$ws.Range("Q4").Value = "D"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("C5").Value = "C:\Users\husmith\Documents\GitHub\Action Pack\Quiz\Question3.png"
$ws.Range("D5").Value = "Text"

# --- Row 6 (new "Text & Image" / SelectMulti question) --------------------
$ws.Range("A6").Value = "Text & Image"
$ws.Range("B6").Value = "Whats the Deal?"
$ws.Range("D6").Value = "SelectMulti"
$ws.Range("K6").Value = "What?"
$ws.Range("L6").Value = $true
$ws.Range("M6").Value = "Who"
$ws.Range("N6").Value = $true
$ws.Range("O6").Value = "Potato"
$ws.Range("Q6").Value = "Ship"

# --- Row 7 (new "Text" / True-False question) ------------------------------
$ws.Range("A7").Value = "Text"
$ws.Range("B7").Value = "Holla at ya boy!"
$ws.Range("D7").Value = "True/False"
$ws.Range("U7").Value = "Flase"

# --- Row heights: rows 3-5 grow to fit the longer picture-path text -------
$ws.Rows.Item(3).RowHeight = 75
$ws.Rows.Item(4).RowHeight = 75
$ws.Rows.Item(5).RowHeight = 75

# --- Update the active selection/view -------------------------------------
$ws.Range("V5").Select() | Out-Null
